$wb = $excel.ActiveWorkbook

# Sheet "展览" - update column F (想去人数) values for rows 3-6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 250
$ws1.Range("F4").Value = 2642
$ws1.Range("F5").Value = 46
$ws1.Range("F6").Value = 564

# Sheet "全部类型" - update column F (想去人数) values for rows 5-8
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 250
$ws4.Range("F6").Value = 2642
$ws4.Range("F7").Value = 46
$ws4.Range("F8").Value = 564
